# Remove the second argument ("test") from the HYPERLINK() formulas
# in columns S, T, V, W, X, Y for rows 2-4, leaving only the URL argument.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("S", "T", "V", "W", "X", "Y")
$rows = @(2, 3, 4)

foreach ($row in $rows) {
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$row")
        $formula = $cell.Formula
        if ($formula) {
            $newFormula = $formula -replace ';\s*"test"\)', ')'
            $cell.Formula = $newFormula
        }
    }
}
